$wb = $excel.ActiveWorkbook

# --- Rename sheets: Sheet3 -> lipids, Sheet4 -> files ---
$ws1 = $wb.Worksheets.Item("Sheet3")
$ws2 = $wb.Worksheets.Item("Sheet4")
$ws1.Name = "lipids"
$ws2.Name = "files"

# --- Clear the "Theoretical mass" column (J) on the lipids sheet ---
# J1 is the header ("Theoretical mass") and J2:J39 hold computed values.
# ClearContents keeps the cell (with its style) when one was already
# assigned, and drops the cell entirely when it had no style - matching
# the target workbook exactly.
$ws1.Range("J1:J39").ClearContents()

# --- Update selection on the lipids sheet to match the column-J click ---
$null = $ws1.Range("J1:J1048576").Select()
